# AgileJava: Lesson 11: Samples
# Fill in the "实际" (Actual) page-count values for chapter 11/12/13 (F4:F6),
# add a new "G" column that computes the page delta between consecutive
# chapters (mirrors existing column E, which does the same for column D),
# move the selection to F6, and reposition/resize the embedded chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the actual page counts that were previously placeholders (0) ---
$ws.Range("F4").Value = 295
$ws.Range("F5").Value = 332
$ws.Range("F6").Value = 332

# --- New column G: amount of pages read since the previous milestone ---
$ws.Range("G2").Formula = "=F3-F2"
$ws.Range("G3").Formula = "=F4-F3"
$ws.Range("G4").Formula = "=F5-F4"
$ws.Range("G5").Formula = "=F6-F5"
$ws.Range("G6").Value = 0
$ws.Range("G7").Formula = "=F8-F7"

# --- Move the active selection to F6 ---
$ws.Range("F6").Select()

# --- Reposition / resize the embedded chart to its new anchor ---
# from: column I (idx 8), 61911 EMU offset, row 1 (idx 0), 0 EMU offset
# to:   column R (idx 17), 619124 EMU offset, row 23 (idx 22), 171450 EMU offset
$EMU_PER_POINT = 12700

$fromCell = $ws.Cells.Item(0 + 1, 8 + 1)
$toCell = $ws.Cells.Item(22 + 1, 17 + 1)

$chartLeft = $fromCell.Left + (61911 / $EMU_PER_POINT)
$chartTop = $fromCell.Top + (0 / $EMU_PER_POINT)
$chartRight = $toCell.Left + (619124 / $EMU_PER_POINT)
$chartBottom = $toCell.Top + (171450 / $EMU_PER_POINT)

$co = $ws.ChartObjects().Item(1)
$co.Left = $chartLeft
$co.Top = $chartTop
$co.Width = $chartRight - $chartLeft
$co.Height = $chartBottom - $chartTop
